$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

function New-DataSheet($afterSheet) {
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Range("A1").Value = "TestCaseName"
    $ws.Range("B1").Value = "Name"
    $ws.Range("C1").Value = "Email"
    $ws.Range("D1").Value = "Current Address"
    $ws.Range("E1").Value = "Permanent Address"
    $ws.Range("A2").Value = "TC01"
    $ws.Range("B2").Value = "ramesh"
    $ws.Range("C2").Value = "ramesh@gmail.com"
    $ws.Range("D2").Value = "TamilNadu, India"
    $ws.Range("E2").Value = "TamilNadu"
    $ws.Range("A3").Value = "TC02"
    $ws.Range("B3").Value = "vikram"
    $ws.Range("C3").Value = "vikram@gmail.com"
    $ws.Range("D3").Value = "Pune,India"
    $ws.Range("E3").Value = "Pune"

    $ws1.Range("A1:E1").Copy()
    $ws.Range("A1:E1").PasteSpecial(-4122)
    $ws1.Range("A2:E2").Copy()
    $ws.Range("A2:E3").PasteSpecial(-4122)

    $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ramesh@gmail.com")
    $ws.Hyperlinks.Add($ws.Range("C3"), "mailto:vikram@gmail.com")
    $ws1.Range("C2").Copy()
    $ws.Range("C2:C3").PasteSpecial(-4122)

    return $ws
}

$ws2 = New-DataSheet($ws1)
$ws2.Name = "FormsPage"

# copy ws2 (already fully corrected) for the remaining sheets to avoid extra dead styles
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "AlertsFramesWindowsPage"

$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "WidgetsPage"

$ws4.Copy($null, $ws4)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "InteractionsPage"

for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
  Write-Host $i ":" $wb.Worksheets.Item($i).Name
}
